# Add a "Test Status" column (AA) showing FAIL (highlighted) for every
# data row, so a failing/erroring "select all organization" run is obvious.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA1").Value = "Test Status"

$lastRow = 5
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 27)
    $cell.Value = "FAIL"
    $cell.Interior.ColorIndex = 3
}
